$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the data currently in column C (C1:C16) into column B (B1:B16),
# then clear the old column C contents.
for ($r = 1; $r -le 16; $r++) {
    $val = $ws.Cells.Item($r, 3).Value2
    $ws.Cells.Item($r, 2).Value2 = $val
}
$ws.Range("C1:C16").ClearContents()

# Update the selected cell to match the new active selection.
$ws.Range("D3").Select()
